# Apply the "first full iteration finished" edit:
#  1. Bump the cached date-placeholder text (15.03.2015 -> 29.03.2015)
#     on the slide master, every slide layout and the notes master.
#  2. Resize/retext the "0*" textbox into an "update" textbox and move it.
#  3. Remove the "1*" textbox.
#  4. Reshape/move the two existing Freeform curly braces.
#  5. Add a third Freeform curly brace and a second "update" textbox.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$text)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {}
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$newDate = "29.03.2015"

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout off the (single) master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes $newDate
}

# Notes master
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# ---- Slide 1 shape edits -------------------------------------------------
$s = $p.Slides.Item(1)

# TextBox 25 ("0*" -> "update"), resized/moved
$tb25 = $s.Shapes.Item("TextBox 25")
$tb25.Left = 496405 / 12700.0
$tb25.Top = 2954780 / 12700.0
$tb25.Width = 867545 / 12700.0
$tb25.Height = 400110 / 12700.0
$tb25.TextFrame.TextRange.Text = "update"

# TextBox 53 ("1*") removed entirely
$s.Shapes.Item("TextBox 53").Delete()

# Freeform 54 reshaped/moved, flipped vertically
$ff54 = $s.Shapes.Item("Freeform 54")
$ff54.Left = 569305 / 12700.0
$ff54.Top = 5084349 / 12700.0
$ff54.Width = 8042314 / 12700.0
$ff54.Height = 1296979 / 12700.0
$ff54.VerticalFlip = $true

# Freeform 55 reshaped/moved, flipped vertically
$ff55 = $s.Shapes.Item("Freeform 55")
$ff55.Left = 544084 / 12700.0
$ff55.Top = 5445224 / 12700.0
$ff55.Width = 8042314 / 12700.0
$ff55.Height = 1152128 / 12700.0
$ff55.VerticalFlip = $true

# New Freeform (copy of Freeform 55's style/geometry) - lands on shape id 24
# so it has to be the 16th "fresh" id handed out on this slide. Burn through
# the lower ids with throw-away duplicates first.
for ($k = 1; $k -le 15; $k++) {
    $tmp = $ff55.Duplicate()
    $tmp.Delete()
}
$ff23 = $ff55.Duplicate()
$ff23.Name = "Freeform 23"
$ff23.Left = 478835 / 12700.0
$ff23.Top = 4815772 / 12700.0
$ff23.Width = 8042314 / 12700.0
$ff23.Height = 1258903 / 12700.0
$ff23.VerticalFlip = $true

# New TextBox (copy of TextBox 25's style) - lands on shape id 30, i.e. the
# 20th fresh id handed out on this slide.
for ($k = 1; $k -le 3; $k++) {
    $tmp = $ff55.Duplicate()
    $tmp.Delete()
}
$tb29 = $tb25.Duplicate()
$tb29.Name = "TextBox 29"
$tb29.Left = 2256469 / 12700.0
$tb29.Top = 2954780 / 12700.0
$tb29.Width = 867545 / 12700.0
$tb29.Height = 400110 / 12700.0
$tb29.TextFrame.TextRange.Text = "update"
